$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = [double]"0.09494372744984057"
$ws.Range("B3").Value = [double]"0.00295370374963417"
$ws.Range("C3").Value = [double]"0.0007024846405599131"
$ws.Range("D3").Value = [double]"1.386936168491154"
$ws.Range("E3").Value = [double]"0.06648140927507246"
$ws.Range("F3").Value = [double]"0.001576855029932104"
$ws.Range("G3").Value = [double]"0.004330552469336235"
$ws.Range("H3").Value = [double]"0.09789743119947475"
$ws.Range("B4").Value = [double]"0.005014842560012558"
$ws.Range("C4").Value = [double]"0.0009883167805282823"
$ws.Range("D4").Value = [double]"6.026337984629561"
$ws.Range("E4").Value = [double]"0.07144091838492561"
$ws.Range("F4").Value = [double]"0.003077771364769242"
$ws.Range("G4").Value = [double]"0.006951913755255874"
$ws.Range("H4").Value = [double]"0.09995857000985313"
$ws.Range("B5").Value = [double]"0.02968309568782823"
$ws.Range("C5").Value = [double]"0.0052844943183773"
$ws.Range("D5").Value = [double]"8.249544288380871"
$ws.Range("E5").Value = [double]"0.07353830048767952"
$ws.Range("F5").Value = [double]"0.01932564432960181"
$ws.Range("G5").Value = [double]"0.04004054704605465"
$ws.Range("H5").Value = [double]"0.1246268231376688"
$ws.Range("B6").Value = [double]"0.01411365635009341"
$ws.Range("C6").Value = [double]"0.004062603941583091"
$ws.Range("D6").Value = [double]"5.474285235342705"
$ws.Range("E6").Value = [double]"0.1215206144636516"
$ws.Range("F6").Value = [double]"0.006151072678357638"
$ws.Range("G6").Value = [double]"0.02207624002182917"
$ws.Range("H6").Value = [double]"0.109057383799934"
$ws.Range("B7").Value = [double]"0.0137524206252949"
$ws.Range("C7").Value = [double]"0.007369592514327578"
$ws.Range("D7").Value = [double]"5.91177877972546"
$ws.Range("E7").Value = [double]"0.09175672192873137"
$ws.Range("F7").Value = [double]"-0.0006917566528382076"
$ws.Range("G7").Value = [double]"0.02819659790342801"
$ws.Range("H7").Value = [double]"0.1086961480751355"
$ws.Range("B8").Value = [double]"0.01504664014136512"
$ws.Range("C8").Value = [double]"0.005026885048022759"
$ws.Range("D8").Value = [double]"5.014580624175117"
$ws.Range("E8").Value = [double]"0.0755151328422265"
$ws.Range("F8").Value = [double]"0.005194101883118216"
$ws.Range("G8").Value = [double]"0.02489917839961202"
$ws.Range("H8").Value = [double]"0.1099903675912057"
$ws.Range("B9").Value = [double]"0.01411306213272985"
$ws.Range("C9").Value = [double]"0.001729766504107046"
$ws.Range("D9").Value = [double]"4.171840469680987"
$ws.Range("E9").Value = [double]"0.06008115026682132"
$ws.Range("F9").Value = [double]"0.01072277229746749"
$ws.Range("G9").Value = [double]"0.01750335196799221"
$ws.Range("H9").Value = [double]"0.1090567895825704"
$ws.Range("B10").Value = [double]"-0.09494372744984057"
$ws.Range("C10").Value = [double]"0.000466672812801438"
$ws.Range("D10").Value = [double]"-217.8179362130637"
$ws.Range("E10").Value = [double]"0"
$ws.Range("F10").Value = [double]"-0.09585839214829744"
$ws.Range("G10").Value = [double]"-0.09402906275138369"
$ws.Range("B11").Value = [double]"-0.04295641648383502"
$ws.Range("C11").Value = [double]"0.0005041597258985286"
$ws.Range("D11").Value = [double]"-88.69414230081276"
$ws.Range("E11").Value = [double]"0"
$ws.Range("F11").Value = [double]"-0.0439445544186691"
$ws.Range("G11").Value = [double]"-0.04196827854900095"
$ws.Range("H11").Value = [double]"0.05198731096600556"
$ws.Range("B12").Value = [double]"-0.03370071030116659"
$ws.Range("C12").Value = [double]"0.0004911045683600219"
$ws.Range("D12").Value = [double]"-71.94872139639322"
$ws.Range("E12").Value = [double]"7.503291435639763e-243"
$ws.Range("F12").Value = [double]"-0.03466326051081057"
$ws.Range("G12").Value = [double]"-0.0327381600915226"
$ws.Range("H12").Value = [double]"0.06124301714867399"
$ws.Range("B13").Value = [double]"-0.03086050529900056"
$ws.Range("C13").Value = [double]"0.000493229948543745"
$ws.Range("D13").Value = [double]"-64.65286474036544"
$ws.Range("E13").Value = [double]"2.141645923395036e-82"
$ws.Range("F13").Value = [double]"-0.03182722120287207"
$ws.Range("G13").Value = [double]"-0.02989378939512906"
$ws.Range("H13").Value = [double]"0.06408322215084"
$ws.Range("B14").Value = [double]"-0.02542508167357629"
$ws.Range("C14").Value = [double]"0.0004814589471612722"
$ws.Range("D14").Value = [double]"-55.27048763482643"
$ws.Range("E14").Value = [double]"1.238011841610718e-25"
$ws.Range("F14").Value = [double]"-0.02636872675208696"
$ws.Range("G14").Value = [double]"-0.02448143659506562"
$ws.Range("H14").Value = [double]"0.06951864577626428"
$ws.Range("B15").Value = [double]"-0.02167993985317307"
$ws.Range("C15").Value = [double]"0.0004766972168157197"
$ws.Range("D15").Value = [double]"-47.76031884152425"
$ws.Range("E15").Value = [double]"1.109814012445849e-76"
$ws.Range("F15").Value = [double]"-0.02261425208032905"
$ws.Range("G15").Value = [double]"-0.02074562762601711"
$ws.Range("H15").Value = [double]"0.0732637875966675"
$ws.Range("B16").Value = [double]"-0.02066151675137504"
$ws.Range("C16").Value = [double]"0.000473535197716213"
$ws.Range("D16").Value = [double]"-46.88009775912229"
$ws.Range("E16").Value = [double]"2.502583367376303e-34"
$ws.Range("F16").Value = [double]"-0.0215896315167294"
$ws.Range("G16").Value = [double]"-0.01973340198602067"
$ws.Range("H16").Value = [double]"0.07428221069846552"
$ws.Range("B17").Value = [double]"-0.01947143362493723"
$ws.Range("C17").Value = [double]"0.0004746600227375002"
$ws.Range("D17").Value = [double]"-44.45756753989649"
$ws.Range("E17").Value = [double]"1.461036229870701e-26"
$ws.Range("F17").Value = [double]"-0.02040175300972089"
$ws.Range("G17").Value = [double]"-0.01854111424015355"
$ws.Range("H17").Value = [double]"0.07547229382490334"
$ws.Range("B18").Value = [double]"-0.01793574536179926"
$ws.Range("C18").Value = [double]"0.0004845511315064013"
$ws.Range("D18").Value = [double]"-39.6155085352211"
$ws.Range("E18").Value = [double]"1.631851749212731e-30"
$ws.Range("F18").Value = [double]"-0.01888545102777163"
$ws.Range("G18").Value = [double]"-0.01698603969582689"
$ws.Range("H18").Value = [double]"0.07700798208804131"
$ws.Range("B19").Value = [double]"-0.01421563431398322"
$ws.Range("C19").Value = [double]"0.0004802234045627498"
$ws.Range("D19").Value = [double]"-31.92047868603346"
$ws.Range("E19").Value = [double]"2.094883394258601e-07"
$ws.Range("F19").Value = [double]"-0.01515685776099948"
$ws.Range("G19").Value = [double]"-0.01327441086696696"
$ws.Range("H19").Value = [double]"0.08072809313585735"
$ws.Range("B20").Value = [double]"-0.01094263605973822"
$ws.Range("C20").Value = [double]"0.0004878742183757785"
$ws.Range("D20").Value = [double]"-23.99415343854259"
$ws.Range("E20").Value = [double]"0.03009225165862738"
$ws.Range("F20").Value = [double]"-0.01189885486452596"
$ws.Range("G20").Value = [double]"-0.009986417254950482"
$ws.Range("H20").Value = [double]"0.08400109139010235"
$ws.Range("B21").Value = [double]"-0.008346704467281104"
$ws.Range("C21").Value = [double]"0.0004867309044120779"
$ws.Range("D21").Value = [double]"-17.25885500251031"
$ws.Range("E21").Value = [double]"9.215290122077e-06"
$ws.Range("F21").Value = [double]"-0.009300682396549239"
$ws.Range("G21").Value = [double]"-0.007392726538012973"
$ws.Range("H21").Value = [double]"0.08659702298255947"
$ws.Range("B22").Value = [double]"-0.006397042984283123"
$ws.Range("C22").Value = [double]"0.0004827943746424578"
$ws.Range("D22").Value = [double]"-13.52086962859235"
$ws.Range("E22").Value = [double]"0.01178590591122428"
$ws.Range("F22").Value = [double]"-0.007343305432646333"
$ws.Range("G22").Value = [double]"-0.005450780535919912"
$ws.Range("H22").Value = [double]"0.08854668446555745"
$ws.Range("B23").Value = [double]"-0.005569446015615608"
$ws.Range("C23").Value = [double]"0.0004801046588001474"
$ws.Range("D23").Value = [double]"-12.09027599068665"
$ws.Range("E23").Value = [double]"0.1049664832985976"
$ws.Range("F23").Value = [double]"-0.006510436695921423"
$ws.Range("G23").Value = [double]"-0.004628455335309792"
$ws.Range("H23").Value = [double]"0.08937428143422496"
$ws.Range("B24").Value = [double]"-0.004927489255277791"
$ws.Range("C24").Value = [double]"0.0004842968945785082"
$ws.Range("D24").Value = [double]"-10.85651304669801"
$ws.Range("E24").Value = [double]"0.05488916140314968"
$ws.Range("F24").Value = [double]"-0.005876696588532838"
$ws.Range("G24").Value = [double]"-0.003978281922022745"
$ws.Range("H24").Value = [double]"0.09001623819456278"
$ws.Range("B25").Value = [double]"-0.002885749243222201"
$ws.Range("C25").Value = [double]"0.0004796543975472105"
$ws.Range("D25").Value = [double]"-6.268184391985406"
$ws.Range("E25").Value = [double]"0.05126048001173619"
$ws.Range("F25").Value = [double]"-0.003825857415215332"
$ws.Range("G25").Value = [double]"-0.00194564107122907"
$ws.Range("H25").Value = [double]"0.09205797820661837"
$ws.Range("B26").Value = [double]"0.0199138480834419"
$ws.Range("C26").Value = [double]"0.001053611888110036"
$ws.Range("D26").Value = [double]"14.17739535684032"
$ws.Range("E26").Value = [double]"0.07490973513595382"
$ws.Range("F26").Value = [double]"0.01784880049504293"
$ws.Range("G26").Value = [double]"0.02197889567184086"
$ws.Range("H26").Value = [double]"0.1148575755332825"
